$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header cells D1/E1 swap their content; F1/G1 (Average Sales / Minimum
# Sales) move down into new rows A17/A18, leaving F1:G1 empty.
$ws.Range("D1").Value = "Minimum Sales"
$ws.Range("E1").Value = "Max Sales"
$ws.Range("F1").Value = $null
$ws.Range("G1").Value = $null

# New standalone rows for the totals, placed after a blank row 16.
$ws.Range("A17").Value = "Total Sales"
$ws.Range("A18").Value = "Average Sales"
